$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 39945
$ws.Range("D2").Value = 57709072
$ws.Range("C3").Value = 95463
$ws.Range("D3").Value = 139875268
$ws.Range("C4").Value = 32518
$ws.Range("D4").Value = 48139554
$ws.Range("C5").Value = 9204
$ws.Range("D5").Value = 13674869
$ws.Range("C6").Value = 2187
$ws.Range("D6").Value = 3248648
$ws.Range("C7").Value = 199
$ws.Range("D7").Value = 293593
$ws.Range("C12").Value = 43291
$ws.Range("D12").Value = 58659852
$ws.Range("C13").Value = 10153
$ws.Range("D13").Value = 14669907
$ws.Range("C14").Value = 27053
$ws.Range("D14").Value = 39649723
$ws.Range("C15").Value = 8618
$ws.Range("D15").Value = 12789365
$ws.Range("C16").Value = 2271
$ws.Range("D16").Value = 3374153
$ws.Range("C20").Value = 10657
$ws.Range("D20").Value = 14065857
$ws.Range("C21").Value = 14025
$ws.Range("D21").Value = 20229835
$ws.Range("C22").Value = 32904
$ws.Range("D22").Value = 48257367
$ws.Range("C23").Value = 10604
$ws.Range("D23").Value = 15759232
$ws.Range("C24").Value = 2765
$ws.Range("D24").Value = 4112174
$ws.Range("C25").Value = 566
$ws.Range("D25").Value = 843092
$ws.Range("C27").Value = 12171
$ws.Range("D27").Value = 16213937
$ws.Range("C28").Value = 8109
$ws.Range("D28").Value = 11727701
$ws.Range("C29").Value = 23513
$ws.Range("D29").Value = 34503428
$ws.Range("C30").Value = 8100
$ws.Range("D30").Value = 12042633
$ws.Range("C31").Value = 2049
$ws.Range("D31").Value = 3056699
$ws.Range("C32").Value = 390
$ws.Range("D32").Value = 582415
$ws.Range("C34").Value = 8701
$ws.Range("D34").Value = 11488796
$ws.Range("C35").Value = 3485
$ws.Range("D35").Value = 5029994
$ws.Range("C36").Value = 8268
$ws.Range("D36").Value = 12077662
$ws.Range("C37").Value = 3302
$ws.Range("D37").Value = 4895961
$ws.Range("C38").Value = 852
$ws.Range("D38").Value = 1269055
$ws.Range("C41").Value = 2625
$ws.Range("D41").Value = 3544328
$ws.Range("C42").Value = 18174
$ws.Range("D42").Value = 26249614
$ws.Range("C43").Value = 53328
$ws.Range("D43").Value = 78139631
$ws.Range("C44").Value = 19625
$ws.Range("D44").Value = 29136980
$ws.Range("C45").Value = 5874
$ws.Range("D45").Value = 8742935
$ws.Range("C46").Value = 1322
$ws.Range("D46").Value = 1972644
$ws.Range("C47").Value = 76
$ws.Range("D47").Value = 112015
$ws.Range("C50").Value = 17582
$ws.Range("D50").Value = 23322312
$ws.Range("C51").Value = 2241
$ws.Range("D51").Value = 3251962
$ws.Range("C52").Value = 7545
$ws.Range("D52").Value = 11086912
$ws.Range("C53").Value = 2520
$ws.Range("D53").Value = 3761684
$ws.Range("C54").Value = 795
$ws.Range("D54").Value = 1187415
$ws.Range("C57").Value = 7627
$ws.Range("D57").Value = 10491587
$ws.Range("C58").Value = 1377
$ws.Range("D58").Value = 2545508
$ws.Range("C59").Value = 3378
$ws.Range("D59").Value = 6263463
$ws.Range("C60").Value = 1329
$ws.Range("D60").Value = 2466462
$ws.Range("C61").Value = 448
$ws.Range("D61").Value = 831083
$ws.Range("C62").Value = 154
$ws.Range("D62").Value = 299600
$ws.Range("C64").Value = 2089
$ws.Range("D64").Value = 3567660
$ws.Range("C65").Value = 16280
$ws.Range("D65").Value = 23502456
$ws.Range("C66").Value = 46826
$ws.Range("D66").Value = 68469913
$ws.Range("C67").Value = 16345
$ws.Range("D67").Value = 24285049
$ws.Range("C68").Value = 4773
$ws.Range("D68").Value = 7109288
$ws.Range("C73").Value = 15706
$ws.Range("D73").Value = 20648742
$ws.Range("C74").Value = 57706
$ws.Range("D74").Value = 83914689
$ws.Range("C75").Value = 159686
$ws.Range("D75").Value = 235102433
$ws.Range("C76").Value = 68475
$ws.Range("D76").Value = 102001330
$ws.Range("C77").Value = 22085
$ws.Range("D77").Value = 32999453
$ws.Range("C78").Value = 5369
$ws.Range("D78").Value = 8019602
$ws.Range("C79").Value = 347
$ws.Range("D79").Value = 515670
$ws.Range("C85").Value = 56787
$ws.Range("D85").Value = 76866280
$ws.Range("C86").Value = 4926
$ws.Range("D86").Value = 7138689
$ws.Range("C87").Value = 12210
$ws.Range("D87").Value = 17933935
$ws.Range("C88").Value = 4032
$ws.Range("D88").Value = 6007958
$ws.Range("C89").Value = 1402
$ws.Range("D89").Value = 2094111
$ws.Range("C90").Value = 318
$ws.Range("D90").Value = 473512
$ws.Range("C93").Value = 5722
$ws.Range("D93").Value = 7680595
$ws.Range("C94").Value = 1720
$ws.Range("D94").Value = 2479136
$ws.Range("C95").Value = 5557
$ws.Range("D95").Value = 8187317
$ws.Range("C96").Value = 2031
$ws.Range("D96").Value = 3023308
$ws.Range("C98").Value = 203
$ws.Range("D98").Value = 306613
$ws.Range("C101").Value = 3803
$ws.Range("D101").Value = 5042690
$ws.Range("C102").Value = 803
$ws.Range("D102").Value = 1440991
$ws.Range("C103").Value = 517
$ws.Range("D103").Value = 972892
$ws.Range("C104").Value = 192
$ws.Range("D104").Value = 359289
$ws.Range("C106").Value = 34
$ws.Range("D106").Value = 69000
$ws.Range("C107").Value = 11420
$ws.Range("D107").Value = 16556636
$ws.Range("C108").Value = 30342
$ws.Range("D108").Value = 44550651
$ws.Range("C109").Value = 10178
$ws.Range("D109").Value = 15131445
$ws.Range("C110").Value = 2810
$ws.Range("D110").Value = 4189510
$ws.Range("C111").Value = 531
$ws.Range("D111").Value = 790903
$ws.Range("C114").Value = 10204
$ws.Range("D114").Value = 13449859
$ws.Range("C115").Value = 32116
$ws.Range("D115").Value = 46288948
$ws.Range("C116").Value = 68951
$ws.Range("D116").Value = 100868930
$ws.Range("C117").Value = 22162
$ws.Range("D117").Value = 32925525
$ws.Range("C118").Value = 6313
$ws.Range("D118").Value = 9399214
$ws.Range("C119").Value = 1211
$ws.Range("D119").Value = 1809092
$ws.Range("C124").Value = 26912
$ws.Range("D124").Value = 35883617
$ws.Range("C125").Value = 38205
$ws.Range("D125").Value = 55101080
$ws.Range("C126").Value = 80559
$ws.Range("D126").Value = 117756402
$ws.Range("C127").Value = 24845
$ws.Range("D127").Value = 36873144
$ws.Range("C128").Value = 6697
$ws.Range("D128").Value = 9951616
$ws.Range("C129").Value = 1369
$ws.Range("D129").Value = 2030792
$ws.Range("C133").Value = 33224
$ws.Range("D133").Value = 44053284
$ws.Range("C134").Value = 14013
$ws.Range("D134").Value = 20279388
$ws.Range("C135").Value = 33725
$ws.Range("D135").Value = 49515046
$ws.Range("C136").Value = 11918
$ws.Range("D136").Value = 17707284
$ws.Range("C137").Value = 3130
$ws.Range("D137").Value = 4664875
$ws.Range("C138").Value = 539
$ws.Range("D138").Value = 802490
$ws.Range("C141").Value = 11275
$ws.Range("D141").Value = 15006382
$ws.Range("C142").Value = 37257
$ws.Range("D142").Value = 53802213
$ws.Range("C143").Value = 85718
$ws.Range("D143").Value = 125547056
$ws.Range("C144").Value = 25511
$ws.Range("D144").Value = 37898430
$ws.Range("C145").Value = 6706
$ws.Range("D145").Value = 10003525
$ws.Range("C146").Value = 1540
$ws.Range("D146").Value = 2289302
$ws.Range("C149").Value = 30625
$ws.Range("D149").Value = 41229299
